$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1049.7778
$ws.Range("I4").Value = 1062.25
$ws.Range("J4").Value = 1039.8
$ws.Range("K4").Value = 1062.25
$ws.Range("L4").Value = 1039.8
$ws.Range("M4").Value = -948.25
$ws.Range("N4").Value = -1267.8
$ws.Range("H19").Value = 400
$ws.Range("I19").Value = 400
$ws.Range("K19").Value = 400
$ws.Range("M19").Value = -225
$ws.Range("H48").Value = 19000
$ws.Range("J48").Value = 19000
$ws.Range("L48").Value = 57000
$ws.Range("N48").Value = -57584
$ws.Range("H56").Value = 19000
$ws.Range("J56").Value = 19000
$ws.Range("L56").Value = 57000
$ws.Range("N56").Value = -58068
$ws.Range("H76").Value = 18523518
$ws.Range("I76").Value = 22732324
$ws.Range("K76").Value = 22732324
$ws.Range("M76").Value = -22732009
$ws.Range("H79").Value = 18523518
$ws.Range("I79").Value = 22732324
$ws.Range("K79").Value = 22732324
$ws.Range("M79").Value = -22731232
$ws.Range("H92").Value = 580.875
$ws.Range("I92").Value = 512.44446
$ws.Range("J92").Value = 786.1667
$ws.Range("K92").Value = 512.44446
$ws.Range("L92").Value = 786.1667
$ws.Range("M92").Value = 735.55554
$ws.Range("N92").Value = -3282.1667
$ws.Range("H96").Value = 3462.7144
$ws.Range("I96").Value = 701.6667
$ws.Range("J96").Value = 20029
$ws.Range("K96").Value = 2105.0001
$ws.Range("L96").Value = 60087
$ws.Range("M96").Value = -732.0001000000002
$ws.Range("N96").Value = -62833
$ws.Range("H97").Value = 1345.6
$ws.Range("J97").Value = 1345.6
$ws.Range("L97").Value = 4036.8
$ws.Range("N97").Value = -5028.799999999999
$ws.Range("H100").Value = 2608.1
$ws.Range("I100").Value = 2373.1428
$ws.Range("J100").Value = 3156.3333
$ws.Range("K100").Value = 2373.1428
$ws.Range("L100").Value = 3156.3333
$ws.Range("M100").Value = -1832.1428
$ws.Range("N100").Value = -4238.3333
$ws.Range("H101").Value = 53512.332
$ws.Range("I101").Value = 964.4
$ws.Range("K101").Value = 2893.2
$ws.Range("M101").Value = -1271.2
$ws.Range("H121").Value = 287557.16
$ws.Range("J121").Value = 287557.16
$ws.Range("L121").Value = 862671.48
$ws.Range("N121").Value = -866165.48
$ws.Range("H123").Value = 74101.42999999999
$ws.Range("J123").Value = 74101.42999999999
$ws.Range("L123").Value = 74101.42999999999
$ws.Range("N123").Value = -83901.42999999999
$ws.Range("H132").Value = 1509.58
$ws.Range("I132").Value = 1405.1904
$ws.Range("K132").Value = 4215.5712
$ws.Range("M132").Value = -1685.5712
$ws.Range("H135").Value = 227.11111
$ws.Range("I135").Value = 238.35294
$ws.Range("K135").Value = 2145.17646
$ws.Range("M135").Value = 389.8235400000003
$ws.Range("H138").Value = 2028.5
$ws.Range("I138").Value = 926.4722
$ws.Range("J138").Value = 2973.0952
$ws.Range("K138").Value = 2779.4166
$ws.Range("L138").Value = 8919.285600000001
$ws.Range("M138").Value = 2360.5834
$ws.Range("N138").Value = -19199.2856
$ws.Range("H141").Value = 1026.8
$ws.Range("I141").Value = 1085.6522
$ws.Range("K141").Value = 3256.9566
$ws.Range("M141").Value = 1923.0434

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 873.5
$ws.Range("I2").Value = 831.6667
$ws.Range("K2").Value = 831.6667
$ws.Range("M2").Value = -718.6667
$ws.Range("H7").Value = 36953.855
$ws.Range("J7").Value = 36953.855
$ws.Range("L7").Value = 36953.855
$ws.Range("N7").Value = -37181.855
$ws.Range("H32").Value = 3589.967
$ws.Range("I32").Value = 2987.6924
$ws.Range("J32").Value = 7203.615
$ws.Range("K32").Value = 2987.6924
$ws.Range("L32").Value = 7203.615
$ws.Range("M32").Value = -2700.6924
$ws.Range("N32").Value = -7777.615
$ws.Range("H45").Value = 11365683
$ws.Range("I45").Value = 2026.5
$ws.Range("K45").Value = 2026.5
$ws.Range("M45").Value = -1649.5
$ws.Range("H52").Value = 53543.4
$ws.Range("J52").Value = 53543.4
$ws.Range("L52").Value = 53543.4
$ws.Range("N52").Value = -54179.4
$ws.Range("H61").Value = 44576.176
$ws.Range("I61").Value = 1147.8636
$ws.Range("J61").Value = 999999
$ws.Range("K61").Value = 1147.8636
$ws.Range("L61").Value = 999999
$ws.Range("M61").Value = -935.8635999999999
$ws.Range("N61").Value = -1000423
$ws.Range("H62").Value = 25000
$ws.Range("J62").Value = 25000
$ws.Range("L62").Value = 25000
$ws.Range("N62").Value = -26248
$ws.Range("H65").Value = 25000
$ws.Range("J65").Value = 25000
$ws.Range("L65").Value = 75000
$ws.Range("N65").Value = -81240
$ws.Range("H74").Value = 47774.047
$ws.Range("I74").Value = 60392.06
$ws.Range("J74").Value = 4872.8
$ws.Range("K74").Value = 60392.06
$ws.Range("L74").Value = 4872.8
$ws.Range("M74").Value = -59518.06
$ws.Range("N74").Value = -6620.8
$ws.Range("H77").Value = 47774.047
$ws.Range("I77").Value = 60392.06
$ws.Range("J77").Value = 4872.8
$ws.Range("K77").Value = 301960.3
$ws.Range("L77").Value = 24364
$ws.Range("M77").Value = -297592.3
$ws.Range("N77").Value = -33100
$ws.Range("H97").Value = 1896.3334
$ws.Range("I97").Value = 1884.5
$ws.Range("J97").Value = 1955.5
$ws.Range("K97").Value = 1884.5
$ws.Range("L97").Value = 1955.5
$ws.Range("M97").Value = -1388.5
$ws.Range("N97").Value = -2947.5
$ws.Range("H102").Value = 111119.6
$ws.Range("I102").Value = 112353.89
$ws.Range("K102").Value = 112353.89
$ws.Range("M102").Value = -110731.89
$ws.Range("H104").Value = 30635
$ws.Range("J104").Value = 30635
$ws.Range("L104").Value = 30635
$ws.Range("N104").Value = -37623
$ws.Range("H110").Value = 1003.65515
$ws.Range("I110").Value = 867.5417
$ws.Range("K110").Value = 867.5417
$ws.Range("M110").Value = 1177.4583
$ws.Range("H116").Value = 873.5
$ws.Range("I116").Value = 831.6667
$ws.Range("K116").Value = 831.6667
$ws.Range("M116").Value = 1462.3333
$ws.Range("H117").Value = 82233.5
$ws.Range("J117").Value = 82233.5
$ws.Range("L117").Value = 82233.5
$ws.Range("N117").Value = -91411.5
$ws.Range("H118").Value = 60698
$ws.Range("J118").Value = 60698
$ws.Range("L118").Value = 60698
$ws.Range("N118").Value = -64012
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H132").Value = 1066.3269
$ws.Range("I132").Value = 900.0714
$ws.Range("J132").Value = 1764.6
$ws.Range("K132").Value = 2700.2142
$ws.Range("L132").Value = 5293.799999999999
$ws.Range("M132").Value = -170.2142000000003
$ws.Range("N132").Value = -10353.8
$ws.Range("H136").Value = 44576.176
$ws.Range("I136").Value = 1147.8636
$ws.Range("J136").Value = 999999
$ws.Range("K136").Value = 3443.5908
$ws.Range("L136").Value = 2999997
$ws.Range("M136").Value = -893.5907999999999
$ws.Range("N136").Value = -3005097
$ws.Range("H140").Value = 97928.60000000001
$ws.Range("J140").Value = 97928.60000000001
$ws.Range("L140").Value = 97928.60000000001
$ws.Range("N140").Value = -108288.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 873.5
$ws.Range("I3").Value = 831.6667
$ws.Range("K3").Value = 831.6667
$ws.Range("M3").Value = -717.6667
$ws.Range("H51").Value = 41491
$ws.Range("J51").Value = 41491
$ws.Range("L51").Value = 41491
$ws.Range("N51").Value = -42473
$ws.Range("H52").Value = 99985
$ws.Range("J52").Value = 99985
$ws.Range("L52").Value = 99985
$ws.Range("N52").Value = -100511
$ws.Range("H53").Value = 25995.5
$ws.Range("J53").Value = 25995.5
$ws.Range("L53").Value = 25995.5
$ws.Range("N53").Value = -27143.5
$ws.Range("H114").Value = 78642.60000000001
$ws.Range("J114").Value = 78642.60000000001
$ws.Range("L114").Value = 78642.60000000001
$ws.Range("N114").Value = -87320.60000000001
$ws.Range("H115").Value = 72329
$ws.Range("J115").Value = 94986
$ws.Range("L115").Value = 94986
$ws.Range("N115").Value = -98120
$ws.Range("H116").Value = 34598.668
$ws.Range("J116").Value = 34598.668
$ws.Range("L116").Value = 34598.668
$ws.Range("N116").Value = -43776.668
$ws.Range("H117").Value = 84181.375
$ws.Range("J117").Value = 84181.375
$ws.Range("L117").Value = 84181.375
$ws.Range("N117").Value = -93359.375
$ws.Range("H119").Value = 37246.75
$ws.Range("J119").Value = 37246.75
$ws.Range("L119").Value = 37246.75
$ws.Range("N119").Value = -46922.75
$ws.Range("H121").Value = 99985
$ws.Range("J121").Value = 99985
$ws.Range("L121").Value = 99985
$ws.Range("N121").Value = -103479
$ws.Range("H127").Value = 59993.5
$ws.Range("J127").Value = 59993.5
$ws.Range("L127").Value = 59993.5
$ws.Range("N127").Value = -69913.5
$ws.Range("H132").Value = 35887.11
$ws.Range("J132").Value = 35887.11
$ws.Range("L132").Value = 35887.11
$ws.Range("N132").Value = -46007.11
$ws.Range("H134").Value = 3147.8262
$ws.Range("I134").Value = 959.06665
$ws.Range("J134").Value = 7251.75
$ws.Range("K134").Value = 2877.19995
$ws.Range("L134").Value = 21755.25
$ws.Range("M134").Value = -342.1999500000002
$ws.Range("N134").Value = -26825.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 25463.857
$ws.Range("J9").Value = 25463.857
$ws.Range("L9").Value = 25463.857
$ws.Range("N9").Value = -25799.857
$ws.Range("H22").Value = 1200.8334
$ws.Range("I22").Value = 1276.625
$ws.Range("J22").Value = 1049.25
$ws.Range("K22").Value = 1276.625
$ws.Range("L22").Value = 1049.25
$ws.Range("M22").Value = -926.625
$ws.Range("N22").Value = -1749.25
$ws.Range("H31").Value = 2225.4424
$ws.Range("I31").Value = 1494.8788
$ws.Range("J31").Value = 3494.3157
$ws.Range("K31").Value = 1494.8788
$ws.Range("L31").Value = 3494.3157
$ws.Range("M31").Value = -1199.8788
$ws.Range("N31").Value = -4084.3157
$ws.Range("H34").Value = 2225.4424
$ws.Range("I34").Value = 1494.8788
$ws.Range("J34").Value = 3494.3157
$ws.Range("K34").Value = 1494.8788
$ws.Range("L34").Value = 3494.3157
$ws.Range("M34").Value = -1292.8788
$ws.Range("N34").Value = -3898.3157
$ws.Range("H58").Value = 1214.1428
$ws.Range("I58").Value = 989.4643
$ws.Range("K58").Value = 989.4643
$ws.Range("M58").Value = -786.4643
$ws.Range("H99").Value = 8377721
$ws.Range("I99").Value = 12349571
$ws.Range("J99").Value = 3909390.5
$ws.Range("K99").Value = 12349571
$ws.Range("L99").Value = 3909390.5
$ws.Range("M99").Value = -12348073
$ws.Range("N99").Value = -3912386.5
$ws.Range("H105").Value = 63946.168
$ws.Range("I105").Value = 101880.63
$ws.Range("J105").Value = 4334.857
$ws.Range("K105").Value = 101880.63
$ws.Range("L105").Value = 4334.857
$ws.Range("M105").Value = -100133.63
$ws.Range("N105").Value = -7828.857
$ws.Range("H107").Value = 1196.8
$ws.Range("J107").Value = 1498.4445
$ws.Range("L107").Value = 1498.4445
$ws.Range("N107").Value = -5338.4445
$ws.Range("H108").Value = 79678.3
$ws.Range("J108").Value = 79678.3
$ws.Range("L108").Value = 79678.3
$ws.Range("N108").Value = -87358.3
$ws.Range("H114").Value = 71972.11
$ws.Range("J114").Value = 71972.11
$ws.Range("L114").Value = 71972.11
$ws.Range("N114").Value = -80650.11
$ws.Range("H118").Value = 55056.5
$ws.Range("J118").Value = 55056.5
$ws.Range("L118").Value = 55056.5
$ws.Range("N118").Value = -58370.5
$ws.Range("H126").Value = 8377721
$ws.Range("I126").Value = 12349571
$ws.Range("J126").Value = 3909390.5
$ws.Range("K126").Value = 37048713
$ws.Range("L126").Value = 11728171.5
$ws.Range("M126").Value = -37046243
$ws.Range("N126").Value = -11733111.5
$ws.Range("H132").Value = 1675.3077
$ws.Range("I132").Value = 1424.4736
$ws.Range("J132").Value = 2356.1428
$ws.Range("K132").Value = 4273.4208
$ws.Range("L132").Value = 7068.428400000001
$ws.Range("M132").Value = -1743.4208
$ws.Range("N132").Value = -12128.4284
$ws.Range("H134").Value = 18989.543
$ws.Range("I134").Value = 1912.2766
$ws.Range("K134").Value = 5736.8298
$ws.Range("M134").Value = -3201.8298
$ws.Range("H136").Value = 1214.1428
$ws.Range("I136").Value = 989.4643
$ws.Range("K136").Value = 2968.3929
$ws.Range("M136").Value = -418.3928999999998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 6369.077
$ws.Range("I7").Value = 10068.25
$ws.Range("J7").Value = 450.4
$ws.Range("K7").Value = 30204.75
$ws.Range("L7").Value = 1351.2
$ws.Range("M7").Value = -30092.75
$ws.Range("N7").Value = -1575.2
$ws.Range("H137").Value = 4963.8423
$ws.Range("I137").Value = 3837.6365
$ws.Range("J137").Value = 6512.375
$ws.Range("K137").Value = 11512.9095
$ws.Range("L137").Value = 19537.125
$ws.Range("M137").Value = -6412.9095
$ws.Range("N137").Value = -29737.125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 505000
$ws.Range("I70").Value = 505000
$ws.Range("K70").Value = 505000
$ws.Range("M70").Value = -504730
$ws.Range("H73").Value = 505000
$ws.Range("I73").Value = 505000
$ws.Range("K73").Value = 505000
$ws.Range("M73").Value = -504064
$ws.Range("H80").Value = 55561050
$ws.Range("J80").Value = 7532.6665
$ws.Range("L80").Value = 7532.6665
$ws.Range("N80").Value = -9528.666499999999
$ws.Range("H83").Value = 55561050
$ws.Range("J83").Value = 7532.6665
$ws.Range("L83").Value = 37663.3325
$ws.Range("N83").Value = -47647.3325
$ws.Range("H97").Value = 3399
$ws.Range("I97").Value = 2787.889
$ws.Range("J97").Value = 5232.3335
$ws.Range("K97").Value = 2787.889
$ws.Range("L97").Value = 5232.3335
$ws.Range("M97").Value = -2291.889
$ws.Range("N97").Value = -6224.3335
$ws.Range("H107").Value = 1100.2858
$ws.Range("I107").Value = 940.4
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 940.4
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 979.6
$ws.Range("N107").Value = -5340
$ws.Range("H108").Value = 57946.082
$ws.Range("J108").Value = 57946.082
$ws.Range("L108").Value = 57946.082
$ws.Range("N108").Value = -65626.08199999999
$ws.Range("H109").Value = 42430.77
$ws.Range("J109").Value = 42430.77
$ws.Range("L109").Value = 42430.77
$ws.Range("N109").Value = -44510.77
$ws.Range("H110").Value = 84952.25
$ws.Range("J110").Value = 84952.25
$ws.Range("L110").Value = 84952.25
$ws.Range("N110").Value = -93132.25
$ws.Range("H116").Value = 59997.332
$ws.Range("J116").Value = 59997.332
$ws.Range("L116").Value = 59997.332
$ws.Range("N116").Value = -69175.33199999999
$ws.Range("H122").Value = 60912.24
$ws.Range("I122").Value = 74368.414
$ws.Range("J122").Value = 3723.5
$ws.Range("K122").Value = 223105.242
$ws.Range("L122").Value = 11170.5
$ws.Range("M122").Value = -220655.242
$ws.Range("N122").Value = -16070.5
$ws.Range("H132").Value = 2559.0923
$ws.Range("I132").Value = 1965.9387
$ws.Range("J132").Value = 4375.625
$ws.Range("K132").Value = 5897.8161
$ws.Range("L132").Value = 13126.875
$ws.Range("M132").Value = -3367.8161
$ws.Range("N132").Value = -18186.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1369.7
$ws.Range("I16").Value = 1174.375
$ws.Range("K16").Value = 1174.375
$ws.Range("M16").Value = -1004.375
$ws.Range("H55").Value = 3998.55
$ws.Range("I55").Value = 883.0769
$ws.Range("J55").Value = 9784.429
$ws.Range("K55").Value = 883.0769
$ws.Range("L55").Value = 9784.429
$ws.Range("M55").Value = -710.0769
$ws.Range("N55").Value = -10130.429
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21248
$ws.Range("H64").Value = 15376.637
$ws.Range("J64").Value = 15376.637
$ws.Range("L64").Value = 15376.637
$ws.Range("N64").Value = -15826.637
$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66240
$ws.Range("H67").Value = 15376.637
$ws.Range("J67").Value = 15376.637
$ws.Range("L67").Value = 15376.637
$ws.Range("N67").Value = -16936.637
$ws.Range("H76").Value = 15143.143
$ws.Range("J76").Value = 16621.6
$ws.Range("L76").Value = 16621.6
$ws.Range("N76").Value = -17297.6
$ws.Range("H79").Value = 15143.143
$ws.Range("J79").Value = 16621.6
$ws.Range("L79").Value = 16621.6
$ws.Range("N79").Value = -18961.6
$ws.Range("H117").Value = 32506.6
$ws.Range("J117").Value = 32506.6
$ws.Range("L117").Value = 32506.6
$ws.Range("N117").Value = -41684.6
$ws.Range("H118").Value = 64500.4
$ws.Range("J118").Value = 64500.4
$ws.Range("L118").Value = 64500.4
$ws.Range("N118").Value = -67814.39999999999
$ws.Range("H121").Value = 39549
$ws.Range("J121").Value = 39549
$ws.Range("L121").Value = 39549
$ws.Range("N121").Value = -43043
$ws.Range("H132").Value = 2118.9143
$ws.Range("J132").Value = 3165.6667
$ws.Range("L132").Value = 9497.000100000001
$ws.Range("N132").Value = -14557.0001
$ws.Range("H136").Value = 3252.9048
$ws.Range("I136").Value = 3155.8386
$ws.Range("J136").Value = 3526.4546
$ws.Range("K136").Value = 9467.515800000001
$ws.Range("L136").Value = 10579.3638
$ws.Range("M136").Value = -6917.515800000001
$ws.Range("N136").Value = -15679.3638
$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 50000
$ws.Range("K137").Value = 50000
$ws.Range("M137").Value = -44900

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 13737.111
$ws.Range("J63").Value = 12576.429
$ws.Range("L63").Value = 12576.429
$ws.Range("N63").Value = -13824.429
$ws.Range("H66").Value = 13737.111
$ws.Range("J66").Value = 12576.429
$ws.Range("K66").Value = 53398.5
$ws.Range("L66").Value = 37729.287
$ws.Range("N66").Value = -43969.287
$ws.Range("H81").Value = 1558.6666
$ws.Range("I81").Value = 870
$ws.Range("J81").Value = 5002
$ws.Range("K81").Value = 1740
$ws.Range("L81").Value = 10004
$ws.Range("M81").Value = -679
$ws.Range("N81").Value = -12126
$ws.Range("H82").Value = 47150.5
$ws.Range("J82").Value = 47150.5
$ws.Range("L82").Value = 47150.5
$ws.Range("N82").Value = -47916.5
$ws.Range("H84").Value = 1558.6666
$ws.Range("I84").Value = 870
$ws.Range("J84").Value = 5002
$ws.Range("K84").Value = 8700
$ws.Range("L84").Value = 50020
$ws.Range("M84").Value = -3396
$ws.Range("N84").Value = -60628
$ws.Range("H85").Value = 47150.5
$ws.Range("J85").Value = 47150.5
$ws.Range("L85").Value = 47150.5
$ws.Range("N85").Value = -49802.5
$ws.Range("H113").Value = 2242.2856
$ws.Range("I113").Value = 2739.4
$ws.Range("K113").Value = 8218.200000000001
$ws.Range("M113").Value = -6048.200000000001
$ws.Range("H121").Value = 67161
$ws.Range("J121").Value = 67161
$ws.Range("L121").Value = 67161
$ws.Range("N121").Value = -70655
$ws.Range("H132").Value = 1403463.5
$ws.Range("I132").Value = 878.75
$ws.Range("K132").Value = 2636.25
$ws.Range("M132").Value = -106.25
